# Append four new transaction rows (5-8) to the Transactions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("GK7CQRV441", "5",    "08/04/2015", "IrregularIncome",  "Found on the ground"),
    @("RCE8DFA1PV", "300",  "08/04/2015", "IrregularExpense", "Robbery"),
    @("DZYM649Q2T", "25",   "08/04/2015", "RegularExpense",   "TV + INTERNET"),
    @("CUJUW9OMXR", "3.85", "08/04/2015", "IrregularIncome",  "Toto")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        # Force text interpretation so numeric-looking / date-looking
        # values ("5", "300", "08/04/2015", ...) are stored as plain text,
        # matching the rest of the sheet's cells rather than being
        # auto-converted into numbers or date serials.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
        # Revert to the sheet's default style so the new cells don't pick
        # up a distinct "Text" number-format style from the rest of the data.
        $cell.Style = "Normal"
    }
}
